$wb = $excel.ActiveWorkbook

# --- 1. Content fix: remove hyphens in "Scene type / Tasks" values ---
# "0001-Cold shelf, 0002-Ambient section" -> "0001 Cold shelf, 0002 Ambient section"
$ws1 = $wb.Worksheets.Item("Functional KPIs")
$oldText = "0001-Cold shelf, 0002-Ambient section"
$newText = "0001 Cold shelf, 0002 Ambient section"
$used = $ws1.UsedRange
$lastRow = $used.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws1.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}

# --- 2. Defined name "Validation_List" now points at external ref [2] instead of [1] ---
$wb.Names.Item("Validation_List").RefersTo = "=[2]Set_up!`$A`$90:`$A`$124"

# --- 3. View changes: zoom 100% -> 140% on both sheets, selection reset to A1 ---
$ws2 = $wb.Worksheets.Item("Instructions")
[void]$ws2.Activate()
[void]$ws2.Range("A1").Select()
$excel.ActiveWindow.Zoom = 140

[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
$excel.ActiveWindow.Zoom = 140
